$wb = $excel.ActiveWorkbook

# Sheet ALC, Row 106 (Leve Item ID 19903)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 1866.3334
$ws.Range("I106").Value = 1800
$ws.Range("K106").Value = 1800
$ws.Range("M106").Value = -1169

# Sheet ALC, Row 112 (Leve Item ID 27960)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2023.25
$ws.Range("I112").Value = 1099
$ws.Range("J112").Value = 2063.4348
$ws.Range("K112").Value = 3297
$ws.Range("L112").Value = 6190.3044
$ws.Range("M112").Value = -2189
$ws.Range("N112").Value = -8406.304400000001

# Sheet ALC, Row 131 (Leve Item ID 36108)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 2497.75
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

# Sheet ALC, Row 132 (Leve Item ID 44049)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3210.4092
$ws.Range("I132").Value = 2854.3
$ws.Range("J132").Value = 6771.5
$ws.Range("K132").Value = 8562.900000000001
$ws.Range("L132").Value = 20314.5
$ws.Range("M132").Value = -6032.900000000001
$ws.Range("N132").Value = -25374.5

# Sheet ALC, Row 138 (Leve Item ID 44169)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2179.72
$ws.Range("J138").Value = 2315.2812
$ws.Range("L138").Value = 6945.8436
$ws.Range("N138").Value = -17225.8436

# Sheet ARM, Row 61 (Leve Item ID 43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3858.1538
$ws.Range("I61").Value = 3858.1538
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3858.1538
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3646.1538
$ws.Range("N61").ClearContents()

# Sheet ARM, Row 74 (Leve Item ID 44000)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2179.7273
$ws.Range("I74").Value = 1773.8889
$ws.Range("J74").Value = 4006
$ws.Range("K74").Value = 1773.8889
$ws.Range("L74").Value = 4006
$ws.Range("M74").Value = -899.8888999999999
$ws.Range("N74").Value = -5754

# Sheet ARM, Row 77 (Leve Item ID 44000)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2179.7273
$ws.Range("I77").Value = 1773.8889
$ws.Range("J77").Value = 4006
$ws.Range("K77").Value = 8869.4445
$ws.Range("L77").Value = 20030
$ws.Range("M77").Value = -4501.4445
$ws.Range("N77").Value = -28766

# Sheet ARM, Row 132 (Leve Item ID 43997)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2177.3635
$ws.Range("I132").Value = 2196
$ws.Range("J132").Value = 1991
$ws.Range("K132").Value = 6588
$ws.Range("L132").Value = 5973
$ws.Range("M132").Value = -4058
$ws.Range("N132").Value = -11033

# Sheet ARM, Row 136 (Leve Item ID 43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3858.1538
$ws.Range("I136").Value = 3858.1538
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 11574.4614
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -9024.4614
$ws.Range("N136").ClearContents()

# Sheet BSM, Row 20 (Leve Item ID 14149)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1056.4117
$ws.Range("I20").Value = 1125.8572
$ws.Range("K20").Value = 1125.8572
$ws.Range("M20").Value = -878.8571999999999

# Sheet BSM, Row 35 (Leve Item ID 2350)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 50074
$ws.Range("J35").Value = 50074
$ws.Range("L35").Value = 50074
$ws.Range("N35").Value = -50694

# Sheet BSM, Row 57 (Leve Item ID 43233)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 8000
$ws.Range("I57").Value = 8000
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 8000
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -7280
$ws.Range("N57").ClearContents()

# Sheet BSM, Row 136 (Leve Item ID 43233)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H136").Value = 8000
$ws.Range("I136").Value = 8000
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 8000
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2900
$ws.Range("N136").ClearContents()

# Sheet CRP, Row 58 (Leve Item ID 44021)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3494.2222
$ws.Range("I58").Value = 3249.8333
$ws.Range("J58").Value = 3983
$ws.Range("K58").Value = 3249.8333
$ws.Range("L58").Value = 3983
$ws.Range("M58").Value = -3046.8333
$ws.Range("N58").Value = -4389

# Sheet CRP, Row 94 (Leve Item ID 32934)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 3032.4443
$ws.Range("I94").Value = 3763.8333
$ws.Range("K94").Value = 3763.8333
$ws.Range("M94").Value = -3312.8333

# Sheet CRP, Row 136 (Leve Item ID 44021)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3494.2222
$ws.Range("I136").Value = 3249.8333
$ws.Range("J136").Value = 3983
$ws.Range("K136").Value = 9749.499899999999
$ws.Range("L136").Value = 11949
$ws.Range("M136").Value = -7199.499899999999
$ws.Range("N136").Value = -17049

# Sheet CUL, Row 97 (Leve Item ID 19846)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 815.3333
$ws.Range("I97").Value = 1000.6667
$ws.Range("K97").Value = 3002.0001
$ws.Range("M97").Value = -2506.0001

# Sheet CUL, Row 107 (Leve Item ID 27838)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 771.375
$ws.Range("J107").Value = 827.7778
$ws.Range("L107").Value = 2483.3334
$ws.Range("N107").Value = -6323.3334

# Sheet CUL, Row 129 (Leve Item ID 36054)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 3641.9375
$ws.Range("J129").Value = 4252.385
$ws.Range("L129").Value = 12757.155
$ws.Range("N129").Value = -22757.155

# Sheet CUL, Row 131 (Leve Item ID 36060)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 75753.266
$ws.Range("I131").Value = 93672.414
$ws.Range("J131").Value = 4076.6667
$ws.Range("K131").Value = 281017.242
$ws.Range("L131").Value = 12230.0001
$ws.Range("M131").Value = -275977.242
$ws.Range("N131").Value = -22310.0001

# Sheet CUL, Row 137 (Leve Item ID 44088)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 3049.7144
$ws.Range("I137").Value = 2116.3333
$ws.Range("K137").Value = 6348.999899999999
$ws.Range("M137").Value = -1248.999899999999

# Sheet GSM, Row 102 (Leve Item ID 36169)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3013.5518
$ws.Range("I102").Value = 2795.8235
$ws.Range("K102").Value = 2795.8235
$ws.Range("M102").Value = -1173.8235

# Sheet GSM, Row 126 (Leve Item ID 36184)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3492.5334
$ws.Range("I126").Value = 2865.8
$ws.Range("K126").Value = 8597.400000000001
$ws.Range("M126").Value = -6127.400000000001

# Sheet LTW, Row 25 (Leve Item ID 3547)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 1000
$ws.Range("J25").Value = 800
$ws.Range("L25").Value = 800
$ws.Range("N25").Value = -1260

# Sheet LTW, Row 46 (Leve Item ID 5282)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1343
$ws.Range("I46").Value = 1014.5
$ws.Range("K46").Value = 1014.5
$ws.Range("M46").Value = -826.5

# Sheet LTW, Row 136 (Leve Item ID 44060)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2728.1052
$ws.Range("I136").Value = 1603.5454
$ws.Range("J136").Value = 4274.375
$ws.Range("K136").Value = 4810.6362
$ws.Range("L136").Value = 12823.125
$ws.Range("M136").Value = -2260.6362
$ws.Range("N136").Value = -17923.125

# Sheet WVR, Row 81 (Leve Item ID 12596)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1592.5
$ws.Range("I81").Value = 1344
$ws.Range("J81").Value = 2503.6667
$ws.Range("K81").Value = 2688
$ws.Range("L81").Value = 5007.3334
$ws.Range("M81").Value = -1627
$ws.Range("N81").Value = -7129.3334

# Sheet WVR, Row 84 (Leve Item ID 12596)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1592.5
$ws.Range("I84").Value = 1344
$ws.Range("J84").Value = 2503.6667
$ws.Range("K84").Value = 13440
$ws.Range("L84").Value = 25036.667
$ws.Range("M84").Value = -8136
$ws.Range("N84").Value = -35644.667

# Sheet WVR, Row 123 (Leve Item ID 34127)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 67200
$ws.Range("J123").Value = 67200
$ws.Range("L123").Value = 67200
$ws.Range("N123").Value = -77000

# Sheet WVR, Row 136 (Leve Item ID 44031)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1931.3846
$ws.Range("I136").Value = 1191.7273
$ws.Range("J136").Value = 5999.5
$ws.Range("K136").Value = 3575.1819
$ws.Range("L136").Value = 17998.5
$ws.Range("M136").Value = -1025.1819
$ws.Range("N136").Value = -23098.5
